$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")

$ws.Range("B14").Value = "cus_N9vHRyiv9GsnhO"
$ws.Range("B18").Value = "cus_N9vHDHWtVlR36c"
